$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 26 de Agosto de 2020 a las 23:13"

# --- Costa Rica / Azerbaiyan swapped position in the country list ---
# Row 64 used to be Azerbaiyan, row 65 used to be Costa Rica.
# Costa Rica now sorts before Azerbaiyan, so row 64 becomes Costa Rica
# (with refreshed case numbers) and row 65 becomes Azerbaiyan (unchanged
# numbers, now sitting one row further down).
$ws.Range("A64").Value = "Costa Rica"
$ws.Range("A65").Value = "Azerbaiyan"

# --- Updated COVID-19 statistics per country (columns B..H) ---

# Estados Unidos (row 4)
$ws.Range("B4").Value = 5991531
$ws.Range("C4").Value = 35803
$ws.Range("D4").Value = 3272292
$ws.Range("E4").Value = 2535922
$ws.Range("G4").Value = 913
$ws.Range("H4").Value = 183317

# Brasil (row 5)
$ws.Range("B5").Value = 3717156
$ws.Range("C5").Value = 42980
$ws.Range("E5").Value = 751096
$ws.Range("G5").Value = 999
$ws.Range("H5").Value = 117665

# Alemania (row 23)
$ws.Range("B23").Value = 239000
$ws.Range("C23").Value = 1428
$ws.Range("E23").Value = 20048
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = 9352

# Israel (row 32)
$ws.Range("B32").Value = 108403
$ws.Range("C32").Value = 1943
$ws.Range("D32").Value = 86466
$ws.Range("E32").Value = 21062

# Belgica (row 40)
$ws.Range("E40").Value = 54278
$ws.Range("H40").Value = 9878

# Guatemala (row 44)
$ws.Range("B44").Value = 70714
$ws.Range("C44").Value = 1063
$ws.Range("D44").Value = 58783
$ws.Range("E44").Value = 9269
$ws.Range("G44").Value = 32
$ws.Range("H44").Value = 2662

# Barein (row 54)
$ws.Range("B54").Value = 50393
$ws.Range("C54").Value = 317
$ws.Range("D54").Value = 47049
$ws.Range("E54").Value = 3158

# Ghana (row 56)
$ws.Range("B56").Value = 43769
$ws.Range("C56").Value = 52
$ws.Range("D56").Value = 42048
$ws.Range("E56").Value = 1451

# Costa Rica (row 64, after reorder) - refreshed numbers
$ws.Range("B64").Value = 36307
$ws.Range("C64").Value = 1002
$ws.Range("D64").Value = 13317
$ws.Range("E64").Value = 22604
$ws.Range("G64").Value = 10
$ws.Range("H64").Value = 386

# Azerbaiyan (row 65, after reorder) - numbers unchanged, just shifted
$ws.Range("B65").Value = 35707
$ws.Range("C65").Value = 148
$ws.Range("D65").Value = 33281
$ws.Range("E65").Value = 1904
$ws.Range("G65").Value = 1
$ws.Range("H65").Value = 522

# Estado de Palestina (row 75)
$ws.Range("E75").Value = 6089
$ws.Range("G75").Value = 4
$ws.Range("H75").Value = 137

# Costa de Marfil (row 79)
$ws.Range("B79").Value = 17603
$ws.Range("C79").Value = 41
$ws.Range("D79").Value = 15941
$ws.Range("E79").Value = 1548

# Mauritania (row 104)
$ws.Range("B104").Value = 6977
$ws.Range("C104").Value = 17
$ws.Range("D104").Value = 6356
$ws.Range("E104").Value = 463

# Zimbabue (row 106)
$ws.Range("B106").Value = 6251
$ws.Range("C106").Value = 55
$ws.Range("D106").Value = 5001
$ws.Range("E106").Value = 1071
$ws.Range("G106").Value = 13
$ws.Range("H106").Value = 179

# Malaui (row 107)
$ws.Range("B107").Value = 5474
$ws.Range("C107").Value = 51
$ws.Range("D107").Value = 3085
$ws.Range("E107").Value = 2216
$ws.Range("G107").Value = 3
$ws.Range("H107").Value = 173

# Republica de Africa Central (row 112)
$ws.Range("B112").Value = 4698
$ws.Range("C112").Value = 7
$ws.Range("D112").Value = 1778
$ws.Range("E112").Value = 2859

# Belice (row 166)
$ws.Range("B166").Value = 760
$ws.Range("C166").Value = 30
$ws.Range("D166").Value = 54
$ws.Range("E166").Value = 695
$ws.Range("G166").Value = 1
$ws.Range("H166").Value = 11
